# "Corregido hasta el capitulo 20"
# Adds new Organizaciones (org20..org42) entries, a couple of newly-known
# names for existing ids, and updates the saved view/selection state so the
# "Organizaciones" sheet is the active tab when the workbook is reopened.

$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("Organizaciones")

# --- New data in the "Organizaciones" sheet ---------------------------------

# Row 21 (id "org20") gained a known name.
$ws3.Cells.Item(21, 2).Value = "publicanos"

# Row 22 (id "org21") gained a known name.
$ws3.Cells.Item(22, 2).Value = "prostitutas"

# Row 23: new id/name pair.
$ws3.Cells.Item(23, 1).Value = "org22"
$ws3.Cells.Item(23, 2).Value = "tribus de Israel"

# Rows 24-43: new ids only (org23 .. org42), no name yet.
$newIds = @("org23","org24","org25","org26","org27","org28","org29","org30", `
            "org31","org32","org33","org34","org35","org36","org37","org38", `
            "org39","org40","org41","org42")

$row = 24
foreach ($id in $newIds) {
    $ws3.Cells.Item($row, 1).Value = $id
    $row = $row + 1
}

# --- View / selection state ---------------------------------------------

# Scroll the "Lugar" sheet's frozen pane further down (selection unchanged).
$ws2 = $wb.Worksheets.Item("Lugar")
$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 41
$ws2.Range("B58").Select()

# Update the remembered selection on "Momento".
$ws4 = $wb.Worksheets.Item("Momento")
$ws4.Activate()
$ws4.Range("A3").Select()

# Finally make "Organizaciones" the active sheet/tab, with its own
# updated selection and scroll position.
$ws3.Activate()
$excel.ActiveWindow.ScrollRow = 7
$ws3.Range("A23").Select()
